# ApprovalLetter.docx merge-field update:
#   1. ${Date}           -> ${currentDate}
#   2. ${Employee Name}  -> ${personFirstName}   (and the "_GoBack" bookmark
#      that used to sit next to "Victor Veteran" now sits right after the
#      new personFirstName merge field)

$d = $word.ActiveDocument

# --- 1) ${Date} -> ${currentDate} -------------------------------------
$rngDate = $d.Content
$rngDate.Find.ClearFormatting()
$rngDate.Find.Replacement.ClearFormatting()
$foundDate = $rngDate.Find.Execute('${Date}', $false, $false, $false, $false, $false, $true, 1, $false, '${currentDate}', 2)
Write-Host "Date field replaced: $foundDate"

# --- 2) Employee Name -> personFirstName, move the _GoBack bookmark ----
$rngName = $d.Content
$rngName.Find.ClearFormatting()
$foundName = $rngName.Find.Execute('Employee Name')
Write-Host "Located 'Employee Name': $foundName"

$oldStart = $rngName.Start
$oldEnd = $rngName.End

# Insert the replacement text right before the old text; inserting at this
# exact boundary makes it inherit the bold/bCs formatting of the preceding
# "${" run, so the new field renders bold like the original did.
$insertPoint = $d.Range($oldStart, $oldStart)
$insertPoint.InsertBefore('personFirstName')
$insertedLen = $insertPoint.End - $insertPoint.Start

# Remove the old "Employee Name" text, now shifted right by the inserted length.
$delRange = $d.Range($oldStart + $insertedLen, $oldEnd + $insertedLen)
$delRange.Delete()

# Re-home the "_GoBack" bookmark immediately after "personFirstName" (it
# previously sat next to "Victor Veteran" further down in the letter).
$bmRange = $d.Range($insertPoint.End, $insertPoint.End)
$d.Bookmarks.Add('_GoBack', $bmRange)
Write-Host "Moved _GoBack bookmark"

Write-Host "Final text: $($d.Content.Text)"
